{"js": "// Append two new paragraphs to the end of the document body:\n//   1. An empty paragraph.\n//   2. A paragraph with the new descriptive text about formatted strings.\n// This mirrors the target OOXML diff, which inserts both paragraphs right\n// before the closing </w:body> (i.e. after the final existing paragraph).\n\nconst newText =\n  \"The use of formatted strings was used in NSLog to define three numeric \" +\n  \"place values: one for the first integer, one for the second integer, \" +\n  \"and the mathematical operation result. By using formatted strings, I \" +\n  \"was able to insert the results quickly into the same string. This \" +\n  \"could have been optimized by using an external function, passing \" +\n  \"associated values, and limiting the lines of code required to generate \" +\n  \"the console statements.\";\n\nconst body = context.document.body;\n\n// Insert the blank paragraph at the very end of the body.\nbody.insertParagraph(\"\", Word.InsertLocation.end);\n\n// Insert the paragraph that will hold the descriptive text (content set\n// below so we can control the exact OOXML, including xml:space=\"preserve\").\nconst textParagraph = body.insertParagraph(newText, Word.InsertLocation.end);\nawait context.sync();\n\n// Rewrite the new paragraph's OOXML directly so the generated markup matches\n// Word's own output exactly: a paragraph with contextualSpacing=0 and a\n// single run (rtl=0) whose <w:t> keeps xml:space=\"preserve\".\nconst escaped = newText\n  .replace(/&/g, \"&amp;\")\n  .replace(/</g, \"&lt;\")\n  .replace(/>/g, \"&gt;\");\n\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr><w:contextualSpacing w:val=\"0\"/></w:pPr>' +\n  '<w:r><w:rPr><w:rtl w:val=\"0\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\">' + escaped + '</w:t>' +\n  '</w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntextParagraph.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Append two new paragraphs to the end of the document:\n#   1. An empty paragraph.\n#   2. A paragraph with the new descriptive text about formatted strings.\n# Matches the target OOXML diff, which inserts both paragraphs right\n# before the closing </w:body> (i.e. after the final existing paragraph).\n\n$d = $word.ActiveDocument\n\n$newText = \"The use of formatted strings was used in NSLog to define three numeric place values: one for the first integer, one for the second integer, and the mathematical operation result. By using formatted strings, I was able to insert the results quickly into the same string. This could have been optimized by using an external function, passing associated values, and limiting the lines of code required to generate the console statements.\"\n\n# Collapsed insertion point at the very end of the document's main story.\n$endPos = $d.Content.End\n$ip = $d.Range($endPos, $endPos)\n\n# Build the replacement OOXML for the two new paragraphs so the run/para\n# properties (contextualSpacing, rtl) and xml:space=\"preserve\" come out\n# exactly as Word itself would author them.\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:contextualSpacing w:val=\"0\"/></w:pPr><w:r><w:rPr><w:rtl w:val=\"0\"/></w:rPr></w:r></w:p><w:p><w:pPr><w:contextualSpacing w:val=\"0\"/></w:pPr><w:r><w:rPr><w:rtl w:val=\"0\"/></w:rPr><w:t xml:space=\"preserve\">' + $newText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$ip.InsertXML($ooxml)\n"}
